# Auto update stock data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update a handful of refreshed EBITDA figures (first row of each company block) ---
# The EBITDA column stores its values as text, so force text formatting while
# writing the new figure, then drop the formatting override again (the source
# cells carry no explicit style) while keeping the value stored as text.
function Set-TextValue($rangeRef, $text) {
    $rng = $ws.Range($rangeRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.ClearFormats()
}

Set-TextValue "B2"  "4.44"
Set-TextValue "B20" "12.04"
Set-TextValue "B32" "25.01"
Set-TextValue "B44" "9.79"
Set-TextValue "B62" "10.70"
Set-TextValue "B68" "11.49"

# --- Add a new "Risk Level" column (J), derived from the Altman Z-Score in column G ---
# Header: copy the formatting of the neighboring header cell, then set the text.
$ws.Range("I1").Copy()
$ws.Range("J1").PasteSpecial(-4122)
$ws.Range("J1").Value = "Risk Level"

# Rows where the Altman Z-Score is below 1.8 -> "High risk"
$ws.Range("J2:J7").Value = "High risk"

# Rows where the Altman Z-Score is between 1.8 and 3.0 -> "Medium risk"
$ws.Range("J8:J13").Value = "Medium risk"
$ws.Range("J26:J31").Value = "Medium risk"
$ws.Range("J32:J37").Value = "Medium risk"
$ws.Range("J68:J73").Value = "Medium risk"

# Rows where the Altman Z-Score is 3.0 or above -> "Low risk"
$ws.Range("J14:J19").Value = "Low risk"
$ws.Range("J20:J25").Value = "Low risk"
$ws.Range("J50:J55").Value = "Low risk"

# Rows with no Altman Z-Score available -> leave the cell present but empty
$blankRows = @(38,39,40,41,42,43,44,45,46,47,48,49,56,57,58,59,60,61,62,63,64,65,66,67,74,75,76,77,78,79)
foreach ($r in $blankRows) {
    $cell = $ws.Cells.Item($r, 10)
    $cell.Font.Bold = $false
}
